$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.691.63"
$ws.Range("E2").Value = "  -1.49%  "

$ws.Range("D3").Value = "1.592.12"
$ws.Range("E3").Value = "  -2.22%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.00%  "

$ws.Range("E6").Value = "  -1.40%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "0.248"
$ws.Range("E8").Value = "  -1.58%  "

$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("D10").Value = "19.66"
$ws.Range("E10").Value = "  -1.70%  "

$ws.Range("D11").Value = "0.0836"
$ws.Range("E11").Value = "  -1.55%  "

$ws.Range("D12").Value = "1.814.93"
$ws.Range("E12").Value = "  -2.27%  "

$ws.Range("D13").Value = "1.588.67"
$ws.Range("E13").Value = "  -1.57%  "

$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -1.41%  "

$ws.Range("E15").Value = "  -1.45%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "63.77"
$ws.Range("E16").Value = "  -1.72%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "26.716.30"
$ws.Range("E17").Value = "  -1.32%  "

$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").Value = "209.65"
$ws.Range("E19").Value = "  -1.60%  "

$ws.Range("E20").Value = "  -0.11%  "

$ws.Range("D21").Value = "6.75"
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("E22").Value = "  -2.27%  "

$ws.Range("D23").Value = "2.36"
$ws.Range("E23").Value = "  -4.69%  "

$ws.Range("E24").Value = "  -2.25%  "

$ws.Range("D25").Value = "146.71"
$ws.Range("E25").Value = "  -0.30%  "

$ws.Range("D26").Value = "7.48"
$ws.Range("E26").Value = "  +2.60%  "

$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("D28").Value = "0.113"
$ws.Range("E28").Value = "  -3.94%  "

$ws.Range("D29").Value = "15.36"
$ws.Range("E29").Value = "  -1.07%  "

$ws.Range("D30").Value = "0.0501"
$ws.Range("E30").Value = "  -0.59%  "

$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  -1.65%  "

$ws.Range("E32").Value = "  -2.71%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "2.94"
$ws.Range("E33").Value = "  -1.43%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.645"
$ws.Range("E34").Value = "  +20.17%  "

$ws.Range("D35").Value = "1.312.50"
$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("E36").Value = "  -3.03%  "

$ws.Range("E37").Value = "  -0.77%  "

$ws.Range("E38").Value = "  -0.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.820"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.22%  "

$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").Value = "0.786"
$ws.Range("E41").Value = "  -2.14%  "

$ws.Range("D42").Value = "2.17"
$ws.Range("E42").Value = "  -3.70%  "

$ws.Range("D43").Value = "5.29"
$ws.Range("E43").Value = "  +0.72%  "

$ws.Range("D44").Value = "63.05"
$ws.Range("E44").Value = "  +1.05%  "

$ws.Range("D45").Value = "1.728.84"
$ws.Range("E45").Value = "  -2.06%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "1.63"
$ws.Range("E46").Value = "  +2.58%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "88.86"
$ws.Range("E47").Value = "  -1.98%  "

$ws.Range("D48").Value = "0.823"
$ws.Range("E48").Value = "  +3.34%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("E49").Value = "  -4.77%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.0508"
$ws.Range("E50").Value = "  -0.94%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0983"
$ws.Range("E51").Value = "  +4.33%  "
